$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-12-05 12:37:58"
}
